$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations")

$ws.Range("A5").Value = "Dog"
$ws.Range("B5").Value = "Dog"
$ws.Range("C5").Value = "Dog"
$ws.Range("D5").Value = "New"

$ws.Range("A6").Value = "Snake"
$ws.Range("B6").Value = "Snake"
$ws.Range("C6").Value = "Snake"
$ws.Range("D6").Value = "New"
